# Auto-generated: update Leviathan Profits price/profit columns (H-N) per sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1042
$ws.Range("I2").Value = 1258.2858
$ws.Range("J2").Value = 285
$ws.Range("K2").Value = 1258.2858
$ws.Range("L2").Value = 285
$ws.Range("M2").Value = -1145.2858
$ws.Range("N2").Value = -511

$ws.Range("H19").Value = 2157.9285
$ws.Range("I19").Value = 789.4545000000001
$ws.Range("J19").Value = 7175.6665
$ws.Range("K19").Value = 789.4545000000001
$ws.Range("L19").Value = 7175.6665
$ws.Range("M19").Value = -614.4545000000001
$ws.Range("N19").Value = -7525.6665

$ws.Range("H121").Value = 5589.6
$ws.Range("J121").Value = 5589.6
$ws.Range("L121").Value = 16768.8
$ws.Range("N121").Value = -20262.8

$ws.Range("H137").Value = 4149.2104
$ws.Range("I137").Value = 3361.1052
$ws.Range("K137").Value = 10083.3156
$ws.Range("M137").Value = -7533.3156

$ws.Range("H138").Value = 1917.7906
$ws.Range("I138").Value = 938.04
$ws.Range("J138").Value = 3278.5557
$ws.Range("K138").Value = 2814.12
$ws.Range("L138").Value = 9835.667099999999
$ws.Range("M138").Value = 2325.88
$ws.Range("N138").Value = -20115.6671


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4488.091
$ws.Range("J2").Value = 3448
$ws.Range("L2").Value = 3448
$ws.Range("N2").Value = -3674

$ws.Range("H5").Value = 202
$ws.Range("J5").Value = 300
$ws.Range("L5").Value = 300
$ws.Range("N5").Value = -524

$ws.Range("H32").Value = 26747.965
$ws.Range("I32").Value = 5331.229
$ws.Range("K32").Value = 5331.229
$ws.Range("M32").Value = -5044.229

$ws.Range("H45").Value = 533958.5
$ws.Range("I45").Value = 843525.75
$ws.Range("J45").Value = 3271.7144
$ws.Range("K45").Value = 843525.75
$ws.Range("L45").Value = 3271.7144
$ws.Range("M45").Value = -843148.75
$ws.Range("N45").Value = -4025.7144

$ws.Range("H74").Value = 1286.6216
$ws.Range("I74").Value = 1102.3334
$ws.Range("K74").Value = 1102.3334
$ws.Range("M74").Value = -228.3334

$ws.Range("H77").Value = 1286.6216
$ws.Range("I77").Value = 1102.3334
$ws.Range("K77").Value = 5511.666999999999
$ws.Range("M77").Value = -1143.666999999999

$ws.Range("H116").Value = 4488.091
$ws.Range("J116").Value = 3448
$ws.Range("L116").Value = 3448
$ws.Range("N116").Value = -8036

$ws.Range("H132").Value = 1829.75
$ws.Range("I132").Value = 1508.1333
$ws.Range("K132").Value = 4524.3999
$ws.Range("M132").Value = -1994.3999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4488.091
$ws.Range("J3").Value = 3448
$ws.Range("L3").Value = 3448
$ws.Range("N3").Value = -3676

$ws.Range("H4").Value = 202
$ws.Range("J4").Value = 300
$ws.Range("L4").Value = 300
$ws.Range("N4").Value = -530

$ws.Range("H20").Value = 11340.941
$ws.Range("I20").Value = 12625.25
$ws.Range("K20").Value = 12625.25
$ws.Range("M20").Value = -12378.25

$ws.Range("H82").Value = 17139.666
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766

$ws.Range("H85").Value = 17139.666
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

$ws.Range("H88").Value = 57483.168
$ws.Range("J88").Value = 57483.168
$ws.Range("L88").Value = 57483.168
$ws.Range("N88").Value = -58295.168

$ws.Range("H91").Value = 57483.168
$ws.Range("J91").Value = 57483.168
$ws.Range("L91").Value = 57483.168
$ws.Range("N91").Value = -60291.168

$ws.Range("H99").Value = 3165.739
$ws.Range("I99").Value = 4770.231
$ws.Range("K99").Value = 4770.231
$ws.Range("M99").Value = -3272.231

$ws.Range("H134").Value = 1436.5
$ws.Range("I134").Value = 1436.5
$ws.Range("K134").Value = 4309.5
$ws.Range("M134").Value = -1774.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37774.215
$ws.Range("I31").Value = 43487.832
$ws.Range("J31").Value = 3492.5
$ws.Range("K31").Value = 43487.832
$ws.Range("L31").Value = 3492.5
$ws.Range("M31").Value = -43192.832
$ws.Range("N31").Value = -4082.5

$ws.Range("H34").Value = 37774.215
$ws.Range("I34").Value = 43487.832
$ws.Range("J34").Value = 3492.5
$ws.Range("K34").Value = 43487.832
$ws.Range("L34").Value = 3492.5
$ws.Range("M34").Value = -43285.832
$ws.Range("N34").Value = -3896.5

$ws.Range("H68").Value = 24618.75
$ws.Range("J68").Value = 24992.857
$ws.Range("L68").Value = 24992.857
$ws.Range("N68").Value = -26490.857

$ws.Range("H71").Value = 24618.75
$ws.Range("J71").Value = 24992.857
$ws.Range("L71").Value = 74978.571
$ws.Range("N71").Value = -82466.571

$ws.Range("H132").Value = 3187.4878
$ws.Range("I132").Value = 2949.1614
$ws.Range("K132").Value = 8847.484199999999
$ws.Range("M132").Value = -6317.484199999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 116.46667
$ws.Range("I12").Value = 10.166667
$ws.Range("J12").Value = 187.33333
$ws.Range("K12").Value = 30.500001
$ws.Range("L12").Value = 561.99999
$ws.Range("M12").Value = 142.499999
$ws.Range("N12").Value = -907.99999

$ws.Range("H37").Value = 111158220
$ws.Range("J37").Value = 111158220
$ws.Range("L37").Value = 333474660
$ws.Range("N37").Value = -333474884

$ws.Range("H137").Value = 4671.2085
$ws.Range("J137").Value = 6689.6924
$ws.Range("L137").Value = 20069.0772
$ws.Range("N137").Value = -30269.0772


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1409.8182
$ws.Range("I102").Value = 1441.16
$ws.Range("K102").Value = 1441.16
$ws.Range("M102").Value = 180.8399999999999

$ws.Range("H113").Value = 3426.5715
$ws.Range("I113").Value = 2999.3333
$ws.Range("K113").Value = 2999.3333
$ws.Range("M113").Value = -829.3332999999998


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1215.7
$ws.Range("I22").Value = 913.4
$ws.Range("K22").Value = 913.4
$ws.Range("M22").Value = -618.4

$ws.Range("H27").Value = 1215.7
$ws.Range("I27").Value = 913.4
$ws.Range("K27").Value = 913.4
$ws.Range("M27").Value = -806.4

$ws.Range("H93").Value = 19234.79
$ws.Range("I93").Value = 1988.1428
$ws.Range("J93").Value = 67525.39999999999
$ws.Range("K93").Value = 1988.1428
$ws.Range("L93").Value = 67525.39999999999
$ws.Range("M93").Value = -740.1428000000001
$ws.Range("N93").Value = -70021.39999999999

$ws.Range("H132").Value = 2374.2678
$ws.Range("I132").Value = 1895.3617
$ws.Range("K132").Value = 5686.0851
$ws.Range("M132").Value = -3156.0851

$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws.Range("H136").Value = 2616.4666
$ws.Range("I136").Value = 2092.139
$ws.Range("K136").Value = 6276.417
$ws.Range("M136").Value = -3726.417


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 350984.5
$ws.Range("I14").Value = 696969
$ws.Range("K14").Value = 696969
$ws.Range("M14").Value = -696801

$ws.Range("H54").Value = 36499

$ws.Range("H132").Value = 901772.5
$ws.Range("I132").Value = 3743.465
$ws.Range("K132").Value = 11230.395
$ws.Range("M132").Value = -8700.395

